$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.35259670456812
$ws.Range("C2").Value = 11.05471106320757
$ws.Range("D2").Value = 3.954914517926018
$ws.Range("F2").Value = 19.67778227069487
$ws.Range("G2").Value = 20.97154831804307
$ws.Range("H2").Value = 12.20958598802698
$ws.Range("B3").Value = 14.55024892360874
$ws.Range("C3").Value = 10.42464483441121
$ws.Range("D3").Value = 3.907313108897109
$ws.Range("F3").Value = 19.78596875342462
$ws.Range("G3").Value = 21.12163510619469
$ws.Range("H3").Value = 12.30714116187034
$ws.Range("B4").Value = 14.03484800502957
$ws.Range("C4").Value = 10.01613215880271
$ws.Range("D4").Value = 3.878092114071695
$ws.Range("F4").Value = 19.8644957982199
$ws.Range("G4").Value = 21.23401038554954
$ws.Range("H4").Value = 12.37131111784548
$ws.Range("B5").Value = 13.81930334204596
$ws.Range("C5").Value = 9.844286739314906
$ws.Range("D5").Value = 3.866196879765517
$ws.Range("F5").Value = 19.89950288517872
$ws.Range("G5").Value = 21.28479675982666
$ws.Range("H5").Value = 12.3985304774682
$ws.Range("B6").Value = 13.78318607602989
$ws.Range("C6").Value = 9.815429823606001
$ws.Range("D6").Value = 3.864222763698191
$ws.Range("F6").Value = 19.9054962985752
$ws.Range("G6").Value = 21.29352855041142
$ws.Range("H6").Value = 12.40311470778541
$ws.Range("B7").Value = 14.03196312018985
$ws.Range("C7").Value = 10.01383624102681
$ws.Range("D7").Value = 3.877931625486328
$ws.Range("F7").Value = 19.86495579092224
$ws.Range("G7").Value = 21.23467521918442
$ws.Range("H7").Value = 12.37167388261005
$ws.Range("B8").Value = 15.08078332551884
$ws.Range("C8").Value = 10.84199488258315
$ws.Range("D8").Value = 3.938506755290284
$ws.Range("F8").Value = 19.71255245303171
$ws.Range("G8").Value = 21.01904728188246
$ws.Range("H8").Value = 12.24233400150661
$ws.Range("B9").Value = 16.94961241506683
$ws.Range("C9").Value = 12.29191387894728
$ws.Range("D9").Value = 4.056869303221662
$ws.Range("F9").Value = 19.51121883876821
$ws.Range("G9").Value = 20.7604943286201
$ws.Range("H9").Value = 12.022782784833
$ws.Range("B10").Value = 18.20032718279803
$ws.Range("C10").Value = 13.24881477896925
$ws.Range("D10").Value = 4.142942508119759
$ws.Range("F10").Value = 19.42468635532752
$ws.Range("G10").Value = 20.67546055803807
$ws.Range("H10").Value = 11.88253583045466
$ws.Range("B11").Value = 18.7415953518525
$ws.Range("C11").Value = 13.66033269598546
$ws.Range("D11").Value = 4.181777102192251
$ws.Range("F11").Value = 19.39901398204984
$ws.Range("G11").Value = 20.66041581423739
$ws.Range("H11").Value = 11.82336945108617
$ws.Range("B12").Value = 18.94250615600754
$ws.Range("C12").Value = 13.81273004559808
$ws.Range("D12").Value = 4.196426148077494
$ws.Range("F12").Value = 19.39128725691687
$ws.Range("G12").Value = 20.65817499740229
$ws.Range("H12").Value = 11.80163643759137
$ws.Range("B13").Value = 18.89941781704681
$ws.Range("C13").Value = 13.78006153937591
$ws.Range("D13").Value = 4.193273905467676
$ws.Range("F13").Value = 19.39286224213386
$ws.Range("G13").Value = 20.65850307444238
$ws.Range("H13").Value = 11.80628704817692
$ws.Range("B14").Value = 18.75820605735063
$ws.Range("C14").Value = 13.67293944196205
$ws.Range("D14").Value = 4.182983485827503
$ws.Range("F14").Value = 19.3983381817929
$ws.Range("G14").Value = 20.66016189971151
$ws.Range("H14").Value = 11.82156796082258
$ws.Range("B15").Value = 18.67117966672821
$ws.Range("C15").Value = 13.60687628047452
$ws.Range("D15").Value = 4.176672601834227
$ws.Range("F15").Value = 19.40195286128415
$ws.Range("G15").Value = 20.6616296165139
$ws.Range("H15").Value = 11.83101565295008
$ws.Range("B16").Value = 18.16439017701667
$ws.Range("C16").Value = 13.22144108917642
$ws.Range("D16").Value = 4.140397205892286
$ws.Range("F16").Value = 19.42664190693616
$ws.Range("G16").Value = 20.67692453846371
$ws.Range("H16").Value = 11.88649620858929
$ws.Range("B17").Value = 17.84634175541544
$ws.Range("C17").Value = 12.97888499060563
$ws.Range("D17").Value = 4.118053630357637
$ws.Range("F17").Value = 19.44531403506329
$ws.Range("G17").Value = 20.69240551261463
$ws.Range("H17").Value = 11.92172236691568
$ws.Range("B18").Value = 17.66080589830431
$ws.Range("C18").Value = 12.83713528092515
$ws.Range("D18").Value = 4.105172645954794
$ws.Range("F18").Value = 19.45734110047372
$ws.Range("G18").Value = 20.70353182015133
$ws.Range("H18").Value = 11.94241920927536
$ws.Range("B19").Value = 17.59754212871152
$ws.Range("C19").Value = 12.78875749081507
$ws.Range("D19").Value = 4.100806621161548
$ws.Range("F19").Value = 19.46163352169303
$ws.Range("G19").Value = 20.70767872978571
$ws.Range("H19").Value = 11.949501459332
$ws.Range("B20").Value = 17.8804685212838
$ws.Range("C20").Value = 13.00493724440003
$ws.Range("D20").Value = 4.120435278859007
$ws.Range("F20").Value = 19.44319292132644
$ws.Range("G20").Value = 20.69052711476213
$ws.Range("H20").Value = 11.91792734839854
$ws.Range("B21").Value = 18.79979395995315
$ws.Range("C21").Value = 13.70449713686395
$ws.Range("D21").Value = 4.186007655795527
$ws.Range("F21").Value = 19.3966754337083
$ws.Range("G21").Value = 20.65958045368019
$ws.Range("H21").Value = 11.81706130077661
$ws.Range("B22").Value = 19.37696243117153
$ws.Range("C22").Value = 14.14167136340438
$ws.Range("D22").Value = 4.228527365845408
$ws.Range("F22").Value = 19.37791095292627
$ws.Range("G22").Value = 20.65952020603583
$ws.Range("H22").Value = 11.75505933407446
$ws.Range("B23").Value = 19.0711026514923
$ws.Range("C23").Value = 13.91017968794877
$ws.Range("D23").Value = 4.205867949139418
$ws.Range("F23").Value = 19.38685328318072
$ws.Range("G23").Value = 20.65769099554869
$ws.Range("H23").Value = 11.7877903108704
$ws.Range("B24").Value = 17.86504816968817
$ws.Range("C24").Value = 12.99316618638614
$ws.Range("D24").Value = 4.119358644885473
$ws.Range("F24").Value = 19.44414785362429
$ws.Range("G24").Value = 20.69136940873242
$ws.Range("H24").Value = 11.91964169093571
$ws.Range("B25").Value = 16.46506296346548
$ws.Range("C25").Value = 11.91858075415373
$ws.Range("D25").Value = 4.024960173249627
$ws.Range("F25").Value = 19.5550357398526
$ws.Range("G25").Value = 20.8123309588212
$ws.Range("H25").Value = 12.07850011500087
